$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("A1").Value = "BusinessName"
$ws.Range("B1").Value = "RegistrationNumber"
$ws.Range("C1").Value = "YearsInBusiness"
$ws.Range("D1").Value = "RevenueUSD"

# Update data row
$ws.Range("A2").Value = "Acme Trading Ltd"
$ws.Range("B2").Value = "RC-00123"
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 250000

# Remove the now-unused Employees column entirely
$ws.Range("E1:E2").EntireColumn.Delete()
